$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 129.41176
$ws.Range("I2").Value = 97.5
$ws.Range("J2").Value = 157.77777
$ws.Range("K2").Value = 97.5
$ws.Range("L2").Value = 157.77777
$ws.Range("M2").Value = 15.5
$ws.Range("N2").Value = -383.77777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 625481.7
$ws.Range("I38").Value = 1250045.2
$ws.Range("J38").Value = 918.125
$ws.Range("K38").Value = 3750135.6
$ws.Range("L38").Value = 2754.375
$ws.Range("M38").Value = -3749763.6
$ws.Range("N38").Value = -3498.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 183763.67
$ws.Range("I43").Value = 50290
$ws.Range("J43").Value = 250500.5
$ws.Range("K43").Value = 50290
$ws.Range("L43").Value = 250500.5
$ws.Range("M43").Value = -50221
$ws.Range("N43").Value = -250638.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 962.8
$ws.Range("I58").Value = 109.61539
$ws.Range("J58").Value = 6508.5
$ws.Range("K58").Value = 328.84617
$ws.Range("L58").Value = 19525.5
$ws.Range("M58").Value = -178.84617
$ws.Range("N58").Value = -19825.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 948.44446
$ws.Range("I98").Value = 908
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 908
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 590
$ws.Range("N98").Value = -4996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1763.6279
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1815.5122
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 5446.536599999999
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -7662.536599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 948.44446
$ws.Range("I122").Value = 908
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2724
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -274
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 12821762
$ws.Range("I135").Value = 406.8421
$ws.Range("J135").Value = 25002050
$ws.Range("K135").Value = 3661.5789
$ws.Range("L135").Value = 225018450
$ws.Range("M135").Value = -1126.5789
$ws.Range("N135").Value = -225023520

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23791.875
$ws.Range("J137").Value = 5526.3184
$ws.Range("L137").Value = 16578.9552
$ws.Range("N137").Value = -21678.9552

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1581.37
$ws.Range("I138").Value = 808.86365
$ws.Range("J138").Value = 2188.3394
$ws.Range("K138").Value = 2426.59095
$ws.Range("L138").Value = 6565.0182
$ws.Range("M138").Value = 2713.40905
$ws.Range("N138").Value = -16845.0182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 813.88
$ws.Range("I61").Value = 787.375
$ws.Range("J61").Value = 1450
$ws.Range("K61").Value = 787.375
$ws.Range("L61").Value = 1450
$ws.Range("M61").Value = -575.375
$ws.Range("N61").Value = -1874

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 35542.535
$ws.Range("I74").Value = 40454.46
$ws.Range("K74").Value = 40454.46
$ws.Range("M74").Value = -39580.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 35542.535
$ws.Range("I77").Value = 40454.46
$ws.Range("K77").Value = 202272.3
$ws.Range("M77").Value = -197904.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1651027.2
$ws.Range("I132").Value = 1927483.1
$ws.Range("K132").Value = 5782449.300000001
$ws.Range("M132").Value = -5779919.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 813.88
$ws.Range("I136").Value = 787.375
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 2362.125
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = 187.875
$ws.Range("N136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33235.37
$ws.Range("I134").Value = 1497.08
$ws.Range("J134").Value = 112581.1
$ws.Range("K134").Value = 4491.24
$ws.Range("L134").Value = 337743.3
$ws.Range("M134").Value = -1956.24
$ws.Range("N134").Value = -342813.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9498.959
$ws.Range("I31").Value = 8957.833
$ws.Range("J31").Value = 10353.368
$ws.Range("K31").Value = 8957.833
$ws.Range("L31").Value = 10353.368
$ws.Range("M31").Value = -8662.833
$ws.Range("N31").Value = -10943.368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9498.959
$ws.Range("I34").Value = 8957.833
$ws.Range("J34").Value = 10353.368
$ws.Range("K34").Value = 8957.833
$ws.Range("L34").Value = 10353.368
$ws.Range("M34").Value = -8755.833
$ws.Range("N34").Value = -10757.368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1032.9286
$ws.Range("I58").Value = 617.3158
$ws.Range("J58").Value = 1910.3334
$ws.Range("K58").Value = 617.3158
$ws.Range("L58").Value = 1910.3334
$ws.Range("M58").Value = -414.3158
$ws.Range("N58").Value = -2316.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1585.7142
$ws.Range("I122").Value = 1585.7142
$ws.Range("K122").Value = 4757.142599999999
$ws.Range("M122").Value = -2307.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1386.5
$ws.Range("I132").Value = 1380.6666
$ws.Range("J132").Value = 1404
$ws.Range("K132").Value = 4141.9998
$ws.Range("L132").Value = 4212
$ws.Range("M132").Value = -1611.9998
$ws.Range("N132").Value = -9272

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1014.4474
$ws.Range("I134").Value = 864.8214
$ws.Range("J134").Value = 1433.4
$ws.Range("K134").Value = 2594.4642
$ws.Range("L134").Value = 4300.200000000001
$ws.Range("M134").Value = -59.46420000000035
$ws.Range("N134").Value = -9370.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1032.9286
$ws.Range("I136").Value = 617.3158
$ws.Range("J136").Value = 1910.3334
$ws.Range("K136").Value = 1851.9474
$ws.Range("L136").Value = 5731.0002
$ws.Range("M136").Value = 698.0526
$ws.Range("N136").Value = -10831.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3403.0278
$ws.Range("J5").Value = 9008.333
$ws.Range("L5").Value = 27024.999
$ws.Range("N5").Value = -27248.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 13334825
$ws.Range("I69").Value = 850
$ws.Range("J69").Value = 15386205
$ws.Range("K69").Value = 2550
$ws.Range("L69").Value = 46158615
$ws.Range("M69").Value = -1739
$ws.Range("N69").Value = -46160237

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 13334825
$ws.Range("I72").Value = 850
$ws.Range("J72").Value = 15386205
$ws.Range("K72").Value = 7650
$ws.Range("L72").Value = 138475845
$ws.Range("M72").Value = -3594
$ws.Range("N72").Value = -138483957

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 30529.646
$ws.Range("I129").Value = 538.75
$ws.Range("J129").Value = 39757.617
$ws.Range("K129").Value = 1616.25
$ws.Range("L129").Value = 119272.851
$ws.Range("M129").Value = 3383.75
$ws.Range("N129").Value = -129272.851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 3600
$ws.Range("J130").Value = 3600
$ws.Range("L130").Value = 10800
$ws.Range("N130").Value = -20840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 54348468
$ws.Range("I131").Value = 411.2857
$ws.Range("J131").Value = 78125740
$ws.Range("K131").Value = 1233.8571
$ws.Range("L131").Value = 234377220
$ws.Range("M131").Value = 3806.1429
$ws.Range("N131").Value = -234387300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 3403.0278
$ws.Range("J135").Value = 9008.333
$ws.Range("L135").Value = 81074.997
$ws.Range("N135").Value = -86144.997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 14495206
$ws.Range("I136").Value = 2628.3333
$ws.Range("J136").Value = 19610234
$ws.Range("K136").Value = 7884.999899999999
$ws.Range("L136").Value = 58830702
$ws.Range("M136").Value = -2784.999899999999
$ws.Range("N136").Value = -58840902

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5293375
$ws.Range("I139").Value = 1484.6154
$ws.Range("J139").Value = 13892697
$ws.Range("K139").Value = 4453.8462
$ws.Range("L139").Value = 41678091
$ws.Range("M139").Value = 686.1538
$ws.Range("N139").Value = -41688371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 17246502
$ws.Range("I140").Value = 1715.1538
$ws.Range("J140").Value = 31257890
$ws.Range("K140").Value = 5145.4614
$ws.Range("L140").Value = 93773670
$ws.Range("M140").Value = 34.53859999999986
$ws.Range("N140").Value = -93784030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 39119.15
$ws.Range("I132").Value = 2332
$ws.Range("J132").Value = 73278.64
$ws.Range("K132").Value = 6996
$ws.Range("L132").Value = 219835.92
$ws.Range("M132").Value = -4466
$ws.Range("N132").Value = -224895.92

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2043.28
$ws.Range("I68").Value = 1742.7858
$ws.Range("K68").Value = 1742.7858
$ws.Range("M68").Value = -993.7858000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2043.28
$ws.Range("I71").Value = 1742.7858
$ws.Range("K71").Value = 8713.929
$ws.Range("M71").Value = -4969.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2760.5144
$ws.Range("I122").Value = 2738.0645
$ws.Range("J122").Value = 2934.5
$ws.Range("K122").Value = 8214.1935
$ws.Range("L122").Value = 8803.5
$ws.Range("M122").Value = -5764.193499999999
$ws.Range("N122").Value = -13703.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 357458.78
$ws.Range("I132").Value = 114415.336
$ws.Range("J132").Value = 630882.7
$ws.Range("K132").Value = 343246.008
$ws.Range("L132").Value = 1892648.1
$ws.Range("M132").Value = -340716.008
$ws.Range("N132").Value = -1897708.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 205248.14
$ws.Range("I136").Value = 286449.4
$ws.Range("J136").Value = 2245
$ws.Range("K136").Value = 859348.2000000001
$ws.Range("L136").Value = 6735
$ws.Range("M136").Value = -856798.2000000001
$ws.Range("N136").Value = -11835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 31950
$ws.Range("J46").Value = 31950
$ws.Range("L46").Value = 31950
$ws.Range("N46").Value = -32412

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2921.4792
$ws.Range("I132").Value = 760.6774
$ws.Range("J132").Value = 6861.7646
$ws.Range("K132").Value = 2282.0322
$ws.Range("L132").Value = 20585.2938
$ws.Range("M132").Value = 247.9677999999999
$ws.Range("N132").Value = -25645.2938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 40692.145
$ws.Range("J133").Value = 40692.145
$ws.Range("L133").Value = 40692.145
$ws.Range("N133").Value = -50812.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 31950
$ws.Range("J134").Value = 31950
$ws.Range("L134").Value = 95850
$ws.Range("N134").Value = -100920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 104910.68
$ws.Range("I136").Value = 704.942
$ws.Range("K136").Value = 2114.826
$ws.Range("M136").Value = 435.174
